$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F14").Value = 831
$ws1.Range("F15").Value = 845
$ws1.Range("F22").Value = 2599
$ws1.Range("F23").Value = 739
$ws1.Range("F25").Value = 2006
$ws1.Range("F26").Value = 462
$ws1.Range("F34").Value = 1009

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F14").Value = 831
$ws4.Range("F23").Value = 2599
$ws4.Range("F24").Value = 739
$ws4.Range("F40").Value = 1009
